$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "35.429.90"
$ws.Range("E2").Value = "  +0.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.897.44"
$ws.Range("E3").Value = "  +0.31%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.77"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.694"
$ws.Range("E6").Value = "  +0.66%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.14"
$ws.Range("E8").Value = "  -0.98%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.358"
$ws.Range("E9").Value = "  +2.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "56.17"
$ws.Range("E10").Value = "  +8.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0756"
$ws.Range("E11").Value = "  +1.99%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0982"
$ws.Range("E12").Value = "  +1.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.08"
$ws.Range("E13").Value = "  +7.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.798"
$ws.Range("E14").Value = "  +10.82%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.171.91"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.01"
$ws.Range("E16").Value = "  +1.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.886.49"
$ws.Range("E17").Value = "  +0.01%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "35.498.09"
$ws.Range("E18").Value = "  +0.63%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "73.72"
$ws.Range("E19").Value = "  +0.77%  "
$ws.Range("E20").Value = "  +1.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "244.75"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "12.99"
$ws.Range("E22").Value = "  +1.37%  "
$ws.Range("E23").Value = "  +4.86%  "
$ws.Range("E24").Value = "  +5.73%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.18"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.50"
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.63"
$ws.Range("E28").Value = "  +1.66%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.35"
$ws.Range("E29").Value = "  +0.71%  "
$ws.Range("E30").Value = "  +0.95%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.35"
$ws.Range("E31").Value = "  +2.26%  "
$ws.Range("E32").Value = "  +4.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.25"
$ws.Range("E33").Value = "  +0.86%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.88"
$ws.Range("E34").Value = "  +26.90%  "
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.47"
$ws.Range("E36").Value = "  -16.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.857"
$ws.Range("E37").Value = "  +1.78%  "
$ws.Range("B38").Value = "Kaspa"
$ws.Range("C38").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0745"
$ws.Range("E38").Value = "  +11.34%  "
$ws.Range("B39").Value = "LidoDAOToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.96"
$ws.Range("E39").Value = "  -1.53%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0226"
$ws.Range("E40").Value = "  +6.57%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "99.24"
$ws.Range("E41").Value = "  +0.59%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "17.07"
$ws.Range("E42").Value = "  -1.41%  "
$ws.Range("E43").Value = "  -0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.80"
$ws.Range("E44").Value = "  +15.05%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.327.19"
$ws.Range("E45").Value = "  +2.37%  "
$ws.Range("E46").Value = "  +0.80%  "
$ws.Range("E47").Value = "  -0.66%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.42"
$ws.Range("E48").Value = "  +0.90%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.74"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.40"
$ws.Range("E50").Value = "  +0.14%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "42.62"
$ws.Range("E51").Value = "  -0.73%  "
